# Updates the cryptocurrency price/volume figures on Sheet1 with freshly
# scraped values (coinranking.com), mirroring the GitHub Actions refresh job.
# Each target cell is plain text in the source data (prices like "44.149.54"
# use '.' as a thousands separator, not a decimal point, and percentages
# keep their surrounding padding spaces), so every write forces the Text
# number format and then restores the cell to the workbook's default
# "Normal" style, matching the unstyled inline strings already on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '44.149.54' },
    @{ Cell = 'E2'; Value = '  +1.50%  ' },
    @{ Cell = 'D3'; Value = '2.253.17' },
    @{ Cell = 'E3'; Value = '  +0.58%  ' },
    @{ Cell = 'E4'; Value = '  +0.17%  ' },
    @{ Cell = 'D5'; Value = '272.48' },
    @{ Cell = 'E5'; Value = '  +5.54%  ' },
    @{ Cell = 'D6'; Value = '87.76' },
    @{ Cell = 'E6'; Value = '  +10.64%  ' },
    @{ Cell = 'E7'; Value = '  +0.10%  ' },
    @{ Cell = 'E8'; Value = '  +0.10%  ' },
    @{ Cell = 'E9'; Value = '  +1.75%  ' },
    @{ Cell = 'D10'; Value = '45.55' },
    @{ Cell = 'E10'; Value = '  +5.05%  ' },
    @{ Cell = 'E11'; Value = '  +0.44%  ' },
    @{ Cell = 'D12'; Value = '7.69' },
    @{ Cell = 'E12'; Value = '  +8.13%  ' },
    @{ Cell = 'E13'; Value = '  +1.78%  ' },
    @{ Cell = 'D14'; Value = '2.591.59' },
    @{ Cell = 'E14'; Value = '  +0.94%  ' },
    @{ Cell = 'D15'; Value = '15.04' },
    @{ Cell = 'E15'; Value = '  +2.56%  ' },
    @{ Cell = 'D16'; Value = '2.269.57' },
    @{ Cell = 'E16'; Value = '  +1.70%  ' },
    @{ Cell = 'D17'; Value = '0.796' },
    @{ Cell = 'E17'; Value = '  +0.07%  ' },
    @{ Cell = 'D18'; Value = '44.085.47' },
    @{ Cell = 'E18'; Value = '  +1.63%  ' },
    @{ Cell = 'D19'; Value = '0.0000104' },
    @{ Cell = 'E19'; Value = '  -0.85%  ' },
    @{ Cell = 'E20'; Value = '  -0.50%  ' },
    @{ Cell = 'D21'; Value = '70.54' },
    @{ Cell = 'E21'; Value = '  -1.18%  ' },
    @{ Cell = 'E22'; Value = '  +2.50%  ' },
    @{ Cell = 'D23'; Value = '234.45' },
    @{ Cell = 'E23'; Value = '  +1.07%  ' },
    @{ Cell = 'D24'; Value = '8.89' },
    @{ Cell = 'E24'; Value = '  -4.66%  ' },
    @{ Cell = 'D26'; Value = '2.54' },
    @{ Cell = 'E26'; Value = '  +14.28%  ' },
    @{ Cell = 'D27'; Value = '10.85' },
    @{ Cell = 'E27'; Value = '  +0.05%  ' },
    @{ Cell = 'E28'; Value = '  +6.20%  ' },
    @{ Cell = 'E29'; Value = '  -4.86%  ' },
    @{ Cell = 'E30'; Value = '  +5.10%  ' },
    @{ Cell = 'D31'; Value = '175.14' },
    @{ Cell = 'E31'; Value = '  +1.10%  ' },
    @{ Cell = 'D32'; Value = '20.93' },
    @{ Cell = 'E32'; Value = '  +1.82%  ' },
    @{ Cell = 'D33'; Value = '0.0897' },
    @{ Cell = 'E33'; Value = '  +2.84%  ' },
    @{ Cell = 'D34'; Value = '5.40' },
    @{ Cell = 'E34'; Value = '  +2.67%  ' },
    @{ Cell = 'E35'; Value = '  +1.24%  ' },
    @{ Cell = 'E36'; Value = '  +3.08%  ' },
    @{ Cell = 'E37'; Value = '  -4.73%  ' },
    @{ Cell = 'E38'; Value = '  -2.15%  ' },
    @{ Cell = 'D39'; Value = '3.50' },
    @{ Cell = 'E39'; Value = '  +21.99%  ' },
    @{ Cell = 'D40'; Value = '12.71' },
    @{ Cell = 'E40'; Value = '  -4.03%  ' },
    @{ Cell = 'E41'; Value = '  +3.22%  ' },
    @{ Cell = 'D42'; Value = '64.90' },
    @{ Cell = 'E43'; Value = '  +1.77%  ' },
    @{ Cell = 'D44'; Value = '0.205' },
    @{ Cell = 'E44'; Value = '  +0.38%  ' },
    @{ Cell = 'E45'; Value = '  -0.60%  ' },
    @{ Cell = 'D46'; Value = '0.0990' },
    @{ Cell = 'E46'; Value = '  +0.70%  ' },
    @{ Cell = 'D47'; Value = '100.89' },
    @{ Cell = 'E47'; Value = '  -2.93%  ' },
    @{ Cell = 'D48'; Value = '1.20' },
    @{ Cell = 'E48'; Value = '  +4.58%  ' },
    @{ Cell = 'E49'; Value = '  +1.49%  ' },
    @{ Cell = 'D50'; Value = '0.431' },
    @{ Cell = 'E50'; Value = '  -8.52%  ' },
    @{ Cell = 'E51'; Value = '  +0.67%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
